$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "IPC PO" predictions (column C) for rows 2-51, reflecting refactored
# weight handling in DenseLayer/NeuralNetwork (new forward-pass outputs).
$newC = @{
    2 = 30.16580649882934
    3 = 29.91327755397509
    4 = 29.92425137224613
    5 = 29.97212626215743
    6 = 30.00820204159942
    7 = 30.14504543960305
    8 = 30.20597661290189
    9 = 30.49339099702649
    10 = 30.50756586496083
    11 = 30.54467052631041
    12 = 30.73448260021426
    13 = 30.86336916496587
    14 = 31.06634062460518
    15 = 31.23863354218343
    16 = 31.30623121835359
    17 = 31.14093985533413
    18 = 30.91667236706866
    19 = 31.11689710311405
    20 = 31.68030221163314
    21 = 32.55461467098568
    22 = 32.5785735425166
    23 = 32.77063300299165
    24 = 32.97614095206262
    25 = 33.03574193207213
    26 = 33.25553034813504
    27 = 33.5578811770355
    28 = 33.65613047871585
    29 = 33.75997469900986
    30 = 34.34806229484705
    31 = 34.63708486459595
    32 = 35.70270976732756
    33 = 35.96025172326057
    34 = 36.3612699998228
    35 = 36.87154410593212
    36 = 37.02772238973446
    37 = 37.84576247163672
    38 = 38.55465602347439
    39 = 39.13423377791571
    40 = 39.44887619032057
    41 = 39.72291496082726
    42 = 39.76860165895368
    43 = 39.90242503756722
    44 = 40.19112861176873
    45 = 41.19744794723858
    46 = 41.2758958709417
    47 = 41.71022517673048
    48 = 42.38125921823084
    49 = 43.50843356850199
    50 = 43.92855664015034
    51 = 44.19887121284617
}

foreach ($r in $newC.Keys) {
    $ws.Cells.Item($r, 3).Value = $newC[$r]
}

# Recompute dependent columns: DELTA = IPC PO - IPC RO, DELTA^2 = DELTA^2
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=C$r-B$r"
    $ws.Cells.Item($r, 5).Formula = "=D$r^2"
}

# Recompute TOTAL (row 52) and MSE (row 53) summary cells
$ws.Range("C52").Formula = "=SUM(D2:D51)"
$ws.Range("E52").Formula = "=SUM(E2:E51)"
$ws.Range("E53").Formula = "=AVERAGE(E2:E51)"

# Excel auto-recalculates after the script runs; freeze the computed results
# as plain numeric literals (matching the source workbook convention of storing
# values rather than live formulas).
$calcRange = $ws.Range("D2:E51")
$calcRange.Value2 = $calcRange.Value2
$summaryRange = $ws.Range("C52:E53")
$summaryRange.Value2 = $summaryRange.Value2
